$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALZ Policy Assignments 03CY23")

# Copy formatting of the row that will become row 18 (currently row 18) so the
# newly-inserted row 17 inherits the same style/height as its neighbour, then
# insert a blank row above it (pushes rows 18+ down to 19+).
$ws.Rows.Item(18).Copy() | Out-Null
$ws.Rows.Item(17).Insert()
$excel.CutCopyMode = $false

# Match the authored row height for the new row.
$ws.Rows.Item(17).RowHeight = 72

# Populate the new row - "Resource Group and Resource locations should match".
# Values are written in the same order the new shared strings were introduced.
$description = "In order to improve resilience and reliability, you need to be aware of where resources are deployed. To aid this awareness, ensure that the location of the resource group matches the location of the resources it contains."
$name = "Resource Group and Resource locations should match"
$ghFile = "Audit-ResourceRGLocation.json"
$linkText = "Audit resource location matches resource group location - 0a914e76-4921-4c19-b460-a2d36003525a (azadvertizer.net)"

$ws.Range("F17").Value = $description
$ws.Range("B17").Value = $name
$ws.Range("C17").Value = $name
$ws.Range("H17").Value = $ghFile
$ws.Range("I17").Value = $linkText

$ws.Range("A17").Value = "Intermediate Root"
$ws.Range("D17").Value = "Policy"
$ws.Range("E17").Value = "Built-in"
$ws.Range("G17").Value = "Audit"
$ws.Range("J17").Value = 45274

# Add the AzAdvertizer hyperlink for the new row's link cell.
$ws.Hyperlinks.Add($ws.Range("I17"), "https://www.azadvertizer.net/azpolicyadvertizer/0a914e76-4921-4c19-b460-a2d36003525a.html", $null, $null, $linkText) | Out-Null

# The table grew by one row (was A1:J50) - refresh the AutoFilter range.
$ws.AutoFilterMode = $false
$ws.Range("A1:J51").AutoFilter() | Out-Null

# Keep the workbook-level _FilterDatabase defined name in sync with the new range.
foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "='ALZ Policy Assignments 03CY23'!`$A`$1:`$J`$51"
    }
}
